# Projects-Information.xlsx update
# - Convert the lower-cased "little title" column (B) on Sheet1 into a
#   LOWER() formula driven off column A.
# - Make part of the marketing blurb (the word "available") bold rich text.
# - Duplicate Sheet1 into a new "Electrical" sheet, scoped down to only the
#   Electrical project row, ready to be duplicated again for the other teams.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Sheet1: turn column B into a formula that lower-cases column A -------
$ws1.Range("B2").Formula = "=LOWER(A2)"
$ws1.Range("B3").Formula = "=LOWER(A3)"
$ws1.Range("B4").Formula = "=LOWER(A4)"
$ws1.Range("B5").Formula = "=LOWER(A5)"

# --- Sheet1: bold the word "available" inside the marketing blurb ---------
$marketingText = "The marketing team is required to create all the information we have available on the website. They go out and find sponsors as well as coordonate and manage our social medias."
$ws1.Range("C5").Value = $marketingText

$boldStart = $marketingText.IndexOf("available") + 1
$boldLen = "available".Length
$boldRun = $ws1.Range("C5").Characters($boldStart, $boldLen)
$boldRun.Font.Bold = $true

$tailStart = $boldStart + $boldLen
$tailLen = $marketingText.Length - $tailStart + 1
$tailRun = $ws1.Range("C5").Characters($tailStart, $tailLen)
$tailRun.Font.Name = "Aptos Narrow"
$tailRun.Font.Size = 11

# --- Sheet1: selection housekeeping ---------------------------------------
$ws1.Range("B2:B5").Select()

# --- Duplicate Sheet1 into a new "Electrical" sheet ------------------------
$ws1.Copy($null, $ws1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "Electrical"

$ws2.Range("A1").Value = "Project Title"
$ws2.Range("A3").Clear()
$ws2.Range("A4").Clear()
$ws2.Range("A5").Clear()

$ws2.Range("B11").Select()

# --- Re-activate Sheet1 and restore its selection --------------------------
$ws1.Activate()
$ws1.Range("B2:B5").Select()

$wb.Save()
